$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the effect-size (95% CI) column for the "survivors" rows with the
# resubmission's re-derived ratio-style estimates.
$ws.Range("D4").Value = "1.65 (1.24~2.07)"
$ws.Range("D5").Value = "3.47 (2.34~4.59)"
$ws.Range("D6").Value = "1.1 (0.76~1.42"

# Rows 7 and 8 ("Arterial" / "Venous" Blood Gas test counts) are merged into
# a single "Blood Gas test count (per day)" row.
$ws.Range("A7").Value = "Blood Gas test count (per day)"
$ws.Range("D7").Value = "1.44 (1.27~1.62)"

$ws.Range("A8").Value = "Total IV fluid volumn (1st day)"
$ws.Range("B8").Value = "1593 (+/-1476)"
$ws.Range("C8").Value = "1759 (+/-1833)"
$ws.Range("D8").Value = "166 (-114~447)"
$ws.Range("E8").Value = 0.24

# The old row 9 (previously "Total IV fluid volumn (1st day)") is now
# redundant since its data moved up into row 8; remove it.
$ws.Rows("9").Delete()

$ws.Range("A9").Select()
